$wb = $excel.ActiveWorkbook

# --- 1. Update the "Date" metadata value on the Metadata sheet (row 8, col B) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-30T16:36:55+00:00"

# --- 2. Add a new row (row 8) to the "Elements" sheet describing the new
#        EquipementSpecifique.LieuRealisationOffre element ---
$ws = $wb.Worksheets.Item("Elements")

# Clone the formatting of the last existing data row (row 7) onto the new row 8
# so the new row gets the same cell style (borders / wrap / vertical alignment)
# as every other data row.
$ws.Range("A7:AJ7").Copy()
$ws.Range("A8:AJ8").PasteSpecial(-4122)
$ws.Rows.Item(8).UseStandardHeight = $true

# Fill in the textual content for the new row
$ws.Range("A8").Value = "EquipementSpecifique.LieuRealisationOffre"
$ws.Range("B8").Value = "EquipementSpecifique.LieuRealisationOffre"
$ws.Range("F8").Value = "'1"
$ws.Range("G8").Value = "'1"
$ws.Range("K8").Value = "https://interop.esante.gouv.fr/ig/mos/StructureDefinition/LieuRealisationOffre`n"
$ws.Range("L8").Value = "Lien vers la classe LieuRealisationOffre"
$ws.Range("M8").Value = "Lien vers la classe LieuRealisationOffre"
$ws.Range("AF8").Value = "EquipementSpecifique.LieuRealisationOffre"
$ws.Range("AG8").Value = "'1"
$ws.Range("AH8").Value = "'1"

# --- 3. Widen column K (11) so the new "URL" column fits its content ---
$ws.Columns.Item(11).ColumnWidth = 61
